$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Panama"
$ws.Cells.Item(2, 2).Value = 4314767
$ws.Cells.Item(2, 3).Value = 332679
$ws.Cells.Item(2, 4).Value = 5642
$ws.Cells.Item(2, 5).Value = 313783
$ws.Cells.Item(2, 6).Value = 13254
$ws.Cells.Item(2, 7).Value = 500
$ws.Cells.Item(2, 8).Value = 0.017
$ws.Cells.Item(2, 9).Value = 1.696
$ws.Cells.Item(2, 10).Value = 7710.243
$ws.Cells.Item(2, 11).Value = 3.772

$ws.Cells.Item(3, 1).Value = "Dominican Republic"
$ws.Cells.Item(3, 2).Value = 10847910
$ws.Cells.Item(3, 3).Value = 230563
$ws.Cells.Item(3, 4).Value = 2959
$ws.Cells.Item(3, 5).Value = 178146
$ws.Cells.Item(3, 6).Value = 49458
$ws.Cells.Item(3, 7).Value = 795
$ws.Cells.Item(3, 8).Value = 0.013
$ws.Cells.Item(3, 9).Value = 1.283
$ws.Cells.Item(3, 10).Value = 2125.414
$ws.Cells.Item(3, 11).Value = 1.607

$ws.Cells.Item(4, 1).Value = "Costa Rica"
$ws.Cells.Item(4, 2).Value = 5094118
$ws.Cells.Item(4, 3).Value = 200024
$ws.Cells.Item(4, 4).Value = 2730
$ws.Cells.Item(4, 5).Value = 163334
$ws.Cells.Item(4, 6).Value = 33960
$ws.Cells.Item(4, 7).Value = 837
$ws.Cells.Item(4, 8).Value = 0.014
$ws.Cells.Item(4, 9).Value = 1.365
$ws.Cells.Item(4, 10).Value = 3926.568
$ws.Cells.Item(4, 11).Value = 2.465

$ws.Cells.Item(5, 1).Value = "Guatemala"
$ws.Cells.Item(5, 2).Value = 17915568
$ws.Cells.Item(5, 3).Value = 167383
$ws.Cells.Item(5, 4).Value = 6150
$ws.Cells.Item(5, 5).Value = 154446
$ws.Cells.Item(5, 6).Value = 6787
$ws.Cells.Item(5, 7).Value = 104
$ws.Cells.Item(5, 8).Value = 0.037
$ws.Cells.Item(5, 9).Value = 3.674
$ws.Cells.Item(5, 10).Value = 934.288
$ws.Cells.Item(5, 11).Value = 1.532

$ws.Cells.Item(6, 1).Value = "Honduras"
$ws.Cells.Item(6, 2).Value = 9904607
$ws.Cells.Item(6, 3).Value = 160983
$ws.Cells.Item(6, 4).Value = 3893
$ws.Cells.Item(6, 5).Value = 63346
$ws.Cells.Item(6, 6).Value = 93744
$ws.Cells.Item(6, 7).Value = 812
$ws.Cells.Item(6, 8).Value = 0.024
$ws.Cells.Item(6, 9).Value = 2.418
$ws.Cells.Item(6, 10).Value = 1625.335
$ws.Cells.Item(6, 11).Value = 0.866

$ws.Cells.Item(7, 1).Value = "El Salvador"
$ws.Cells.Item(7, 2).Value = 6486205
$ws.Cells.Item(7, 3).Value = 58023
$ws.Cells.Item(7, 4).Value = 1750
$ws.Cells.Item(7, 5).Value = 52688
$ws.Cells.Item(7, 6).Value = 3585
$ws.Cells.Item(7, 7).Value = 595
$ws.Cells.Item(7, 8).Value = 0.03
$ws.Cells.Item(7, 9).Value = 3.016
$ws.Cells.Item(7, 10).Value = 894.56
$ws.Cells.Item(7, 11).Value = 16.597

$ws.Cells.Item(8, 1).Value = "Cuba"
$ws.Cells.Item(8, 2).Value = 11326616
$ws.Cells.Item(8, 3).Value = 39004
$ws.Cells.Item(8, 4).Value = 269
$ws.Cells.Item(8, 5).Value = 33776
$ws.Cells.Item(8, 6).Value = 4959
$ws.Cells.Item(8, 7).Value = 715
$ws.Cells.Item(8, 8).Value = 0.007
$ws.Cells.Item(8, 9).Value = 0.69
$ws.Cells.Item(8, 10).Value = 344.357
$ws.Cells.Item(8, 11).Value = 14.418

$ws.Cells.Item(9, 1).Value = "Jamaica"
$ws.Cells.Item(9, 2).Value = 2961167
$ws.Cells.Item(9, 3).Value = 19305
$ws.Cells.Item(9, 4).Value = 378
$ws.Cells.Item(9, 5).Value = 12635
$ws.Cells.Item(9, 6).Value = 6292
$ws.Cells.Item(9, 7).Value = 270
$ws.Cells.Item(9, 8).Value = 0.02
$ws.Cells.Item(9, 9).Value = 1.958
$ws.Cells.Item(9, 10).Value = 651.939
$ws.Cells.Item(9, 11).Value = 4.291

$ws.Cells.Item(10, 1).Value = "Haiti"
$ws.Cells.Item(10, 2).Value = 11402528
$ws.Cells.Item(10, 3).Value = 12143
$ws.Cells.Item(10, 4).Value = 247
$ws.Cells.Item(10, 5).Value = 9354
$ws.Cells.Item(10, 6).Value = 2542
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = 0.02
$ws.Cells.Item(10, 9).Value = 2.034
$ws.Cells.Item(10, 10).Value = 106.494
$ws.Cells.Item(10, 11).Value = 0

$ws.Cells.Item(11, 1).Value = "Bahamas"
$ws.Cells.Item(11, 2).Value = 393244
$ws.Cells.Item(11, 3).Value = 8311
$ws.Cells.Item(11, 4).Value = 178
$ws.Cells.Item(11, 5).Value = 6931
$ws.Cells.Item(11, 6).Value = 1202
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = 0.021
$ws.Cells.Item(11, 9).Value = 2.142
$ws.Cells.Item(11, 10).Value = 2113.446
$ws.Cells.Item(11, 11).Value = 0

$ws.Cells.Item(12, 1).Value = "Guyana"
$ws.Cells.Item(12, 2).Value = 786552
$ws.Cells.Item(12, 3).Value = 8232
$ws.Cells.Item(12, 4).Value = 186
$ws.Cells.Item(12, 5).Value = 7399
$ws.Cells.Item(12, 6).Value = 647
$ws.Cells.Item(12, 7).Value = 1
$ws.Cells.Item(12, 8).Value = 0.023
$ws.Cells.Item(12, 9).Value = 2.259
$ws.Cells.Item(12, 10).Value = 1046.593
$ws.Cells.Item(12, 11).Value = 0.155

$ws.Cells.Item(13, 1).Value = "Trinidad and Tobago"
$ws.Cells.Item(13, 2).Value = 1399488
$ws.Cells.Item(13, 3).Value = 7646
$ws.Cells.Item(13, 4).Value = 138
$ws.Cells.Item(13, 5).Value = 7351
$ws.Cells.Item(13, 6).Value = 157
$ws.Cells.Item(13, 7).Value = 4
$ws.Cells.Item(13, 8).Value = 0.018
$ws.Cells.Item(13, 9).Value = 1.805
$ws.Cells.Item(13, 10).Value = 546.343
$ws.Cells.Item(13, 11).Value = 2.548

$ws.Cells.Item(14, 1).Value = "Aruba"
$ws.Cells.Item(14, 2).Value = 106766
$ws.Cells.Item(14, 3).Value = 7438
$ws.Cells.Item(14, 4).Value = 68
$ws.Cells.Item(14, 5).Value = 7098
$ws.Cells.Item(14, 6).Value = 272
$ws.Cells.Item(14, 7).Value = 25
$ws.Cells.Item(14, 8).Value = 0.009
$ws.Cells.Item(14, 9).Value = 0.914
$ws.Cells.Item(14, 10).Value = 6966.637
$ws.Cells.Item(14, 11).Value = 9.191

$ws.Cells.Item(15, 1).Value = "Nicaragua"
$ws.Cells.Item(15, 2).Value = 6624554
$ws.Cells.Item(15, 3).Value = 6347
$ws.Cells.Item(15, 4).Value = 171
$ws.Cells.Item(15, 5).Value = 4225
$ws.Cells.Item(15, 6).Value = 1951
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 0.027
$ws.Cells.Item(15, 9).Value = 2.694
$ws.Cells.Item(15, 10).Value = 95.81
$ws.Cells.Item(15, 11).Value = 0

$ws.Cells.Item(16, 1).Value = "Curacao"
$ws.Cells.Item(16, 2).Value = 164093
$ws.Cells.Item(16, 3).Value = 4652
$ws.Cells.Item(16, 4).Value = 22
$ws.Cells.Item(16, 5).Value = 4571
$ws.Cells.Item(16, 6).Value = 59
$ws.Cells.Item(16, 7).Value = 2
$ws.Cells.Item(16, 8).Value = 0.005
$ws.Cells.Item(16, 9).Value = 0.473
$ws.Cells.Item(16, 10).Value = 2834.978
$ws.Cells.Item(16, 11).Value = 3.39

$ws.Cells.Item(17, 1).Value = "Saint Lucia"
$ws.Cells.Item(17, 2).Value = 183627
$ws.Cells.Item(17, 3).Value = 2519
$ws.Cells.Item(17, 4).Value = 23
$ws.Cells.Item(17, 5).Value = 1778
$ws.Cells.Item(17, 6).Value = 718
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 0.009
$ws.Cells.Item(17, 9).Value = 0.913
$ws.Cells.Item(17, 10).Value = 1371.803
$ws.Cells.Item(17, 11).Value = 0

$ws.Cells.Item(18, 1).Value = "Barbados"
$ws.Cells.Item(18, 2).Value = 287375
$ws.Cells.Item(18, 3).Value = 2268
$ws.Cells.Item(18, 4).Value = 24
$ws.Cells.Item(18, 5).Value = 1639
$ws.Cells.Item(18, 6).Value = 605
$ws.Cells.Item(18, 7).Value = 207
$ws.Cells.Item(18, 8).Value = 0.011
$ws.Cells.Item(18, 9).Value = 1.058
$ws.Cells.Item(18, 10).Value = 789.213
$ws.Cells.Item(18, 11).Value = 34.215

$ws.Cells.Item(19, 1).Value = "St Martin"
$ws.Cells.Item(19, 2).Value = 38666
$ws.Cells.Item(19, 3).Value = 1408
$ws.Cells.Item(19, 4).Value = 12
$ws.Cells.Item(19, 5).Value = 1050
$ws.Cells.Item(19, 6).Value = 346
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 0.009
$ws.Cells.Item(19, 9).Value = 0.852
$ws.Cells.Item(19, 10).Value = 3641.442
$ws.Cells.Item(19, 11).Value = 0

$ws.Cells.Item(20, 1).Value = "Antigua and Barbuda"
$ws.Cells.Item(20, 2).Value = 97929
$ws.Cells.Item(20, 3).Value = 443
$ws.Cells.Item(20, 4).Value = 9
$ws.Cells.Item(20, 5).Value = 205
$ws.Cells.Item(20, 6).Value = 229
$ws.Cells.Item(20, 7).Value = 16
$ws.Cells.Item(20, 8).Value = 0.02
$ws.Cells.Item(20, 9).Value = 2.032
$ws.Cells.Item(20, 10).Value = 452.369
$ws.Cells.Item(20, 11).Value = 6.987

$ws.Cells.Item(21, 1).Value = "Cayman Islands"
$ws.Cells.Item(21, 2).Value = 65757
$ws.Cells.Item(21, 3).Value = 416
$ws.Cells.Item(21, 4).Value = 2
$ws.Cells.Item(21, 5).Value = 378
$ws.Cells.Item(21, 6).Value = 36
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(21, 8).Value = 0.005
$ws.Cells.Item(21, 9).Value = 0.481
$ws.Cells.Item(21, 10).Value = 632.632
$ws.Cells.Item(21, 11).Value = 0
